# Apply the refreshed cryptos snapshot (price/volume columns) produced by
# the scheduled GitHub Actions scrape. Values are written as literal text
# (matching the source sheet, which stores Price/Volume(1h) as inline
# strings rather than numbers) so formatting such as trailing zeros and
# the "+"/"-" sign on percentages is preserved exactly.
#
# Numeric-looking Price values (e.g. "571.59") are entered with a leading
# apostrophe, same as a user typing into Excel, to force text storage and
# stop them being reinterpreted as numbers (which would silently drop
# formatting like the trailing zero in "27.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range('D2').Value = '60.415.84'
$ws.Range('E2').Value = '  -1.69%  '
# Row 3
$ws.Range('D3').Value = '3.381.36'
$ws.Range('E3').Value = '  -2.00%  '
# Row 4
$ws.Range('E4').Value = '  -0.07%  '
# Row 5
$ws.Range('D5').Value = '''571.59'
$ws.Range('E5').Value = '  -1.42%  '
# Row 6
$ws.Range('D6').Value = '''141.23'
$ws.Range('E6').Value = '  -4.65%  '
# Row 7
$ws.Range('E7').Value = '  +0.05%  '
# Row 8
$ws.Range('D8').Value = '3.379.94'
$ws.Range('E8').Value = '  -2.07%  '
# Row 9
$ws.Range('E9').Value = '  +0.20%  '
# Row 10
$ws.Range('E10').Value = '  -4.34%  '
# Row 11
$ws.Range('E11').Value = '  -0.57%  '
# Row 12
$ws.Range('E12').Value = '  +0.60%  '
# Row 13
$ws.Range('D13').Value = '3.960.28'
$ws.Range('E13').Value = '  -2.04%  '
# Row 14
$ws.Range('D14').Value = '''28.21'
$ws.Range('E14').Value = '  +0.58%  '
# Row 15
$ws.Range('D15').Value = '''0.123'
$ws.Range('E15').Value = '  +0.81%  '
# Row 16
$ws.Range('D16').Value = '''0.0000171'
$ws.Range('E16').Value = '  -2.45%  '
# Row 17
$ws.Range('D17').Value = '3.376.69'
$ws.Range('E17').Value = '  -2.12%  '
# Row 18
$ws.Range('D18').Value = '60.550.61'
$ws.Range('E18').Value = '  -1.68%  '
# Row 19
$ws.Range('D19').Value = '''6.27'
$ws.Range('E19').Value = '  -0.94%  '
# Row 20
$ws.Range('D20').Value = '''14.11'
$ws.Range('E20').Value = '  -1.73%  '
# Row 21
$ws.Range('D21').Value = '''9.22'
$ws.Range('E21').Value = '  -2.48%  '
# Row 22
$ws.Range('D22').Value = '''389.01'
$ws.Range('E22').Value = '  +0.51%  '
# Row 23
$ws.Range('E23').Value = '  -1.60%  '
# Row 24
$ws.Range('D24').Value = '''73.53'
$ws.Range('E24').Value = '  +1.12%  '
# Row 25
$ws.Range('E25').Value = '  -0.30%  '
# Row 26
$ws.Range('D26').Value = '''0.0000117'
$ws.Range('E26').Value = '  -4.37%  '
# Row 27
$ws.Range('D27').Value = '3.518.10'
$ws.Range('E27').Value = '  -2.14%  '
# Row 28
$ws.Range('D28').Value = '''0.179'
$ws.Range('E28').Value = '  -0.55%  '
# Row 29
$ws.Range('E29').Value = '  -0.02%  '
# Row 30
$ws.Range('D30').Value = '''7.39'
$ws.Range('E30').Value = '  -5.62%  '
# Row 31
$ws.Range('D31').Value = '''8.07'
$ws.Range('E31').Value = '  -2.32%  '
# Row 32
$ws.Range('D32').Value = '''2.15'
$ws.Range('E32').Value = '  -1.12%  '
# Row 33
$ws.Range('D33').Value = '''1.44'
$ws.Range('E33').Value = '  -6.38%  '
# Row 35
$ws.Range('D35').Value = '''23.75'
$ws.Range('E35').Value = '  -0.96%  '
# Row 36
$ws.Range('D36').Value = '''6.95'
$ws.Range('E36').Value = '  -1.86%  '
# Row 37
$ws.Range('D37').Value = '3.407.45'
$ws.Range('E37').Value = '  -2.03%  '
# Row 38
$ws.Range('D38').Value = '''167.29'
$ws.Range('E38').Value = '  +0.54%  '
# Row 39
$ws.Range('D39').Value = '''4.95'
$ws.Range('E39').Value = '  -5.30%  '
# Row 40
$ws.Range('E40').Value = '  -4.07%  '
# Row 41
$ws.Range('D41').Value = '''0.0777'
$ws.Range('E41').Value = '  -1.58%  '
# Row 42
$ws.Range('D42').Value = '''27.00'
$ws.Range('E42').Value = '  +3.73%  '
# Row 43
$ws.Range('D43').Value = '''0.782'
$ws.Range('E43').Value = '  -1.83%  '
# Row 44
$ws.Range('D44').Value = '''0.998'
$ws.Range('E44').Value = '  -0.22%  '
# Row 45
$ws.Range('D45').Value = '''4.47'
$ws.Range('E45').Value = '  -0.95%  '
# Row 46
$ws.Range('D46').Value = '''1.70'
$ws.Range('E46').Value = '  -1.26%  '
# Row 47
$ws.Range('D47').Value = '''41.37'
$ws.Range('E47').Value = '  -2.18%  '
# Row 48
$ws.Range('D48').Value = '2.533.59'
$ws.Range('E48').Value = '  -2.71%  '
# Row 49
$ws.Range('E49').Value = '  -3.45%  '
# Row 50
$ws.Range('D50').Value = '''6.85'
# Row 51
$ws.Range('D51').Value = '''23.11'
$ws.Range('E51').Value = '  -0.45%  '
